# Generate Report for Handback
# The handback for the "3724d988..." file came back out of sync with en-US,
# so the status text and the handback timestamps need to be refreshed across
# the Overview sheet and each per-locale report sheet.

$wb = $excel.ActiveWorkbook

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Handed back: not in sync with en-US"

# --- Overview sheet: both locale status columns for the affected file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("E1:F3").EntireColumn.AutoFit() | Out-Null

# --- zh-cn sheet: Status column + refreshed handback datetime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-09-09 13:11:52"
$wsZhCn.Range("C1:C3").EntireColumn.AutoFit() | Out-Null

# --- de-de sheet: Status column + refreshed handback datetime ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-09-09 13:12:22"
$wsDeDe.Range("C1:C3").EntireColumn.AutoFit() | Out-Null
